{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (first paragraph of the document, the hidden\n// \"**ID__AFFARS_..._ID**\" marker paragraph):\n//   1. pPr gets a <w:pBdr> with all four edges set to w:space=\"5\"\n//      (no line style / width \u2014 just spacing).\n//   2. <w:ind w:left=\"120\"/> -> <w:ind w:left=\"225\"/>  (6pt -> 11.25pt)\n//   3. The two runs \"**ID__AFFARS_pgi_5337_topic_2__ID**\" + \" \" (a\n//      trailing-space-only run) collapse into a single run whose text\n//      becomes \"**ID__AFFARS_AFICC_PGI_5337__ID**\" (no trailing space),\n//      keeping the original run formatting (font/color/size).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// --- 1) & 3): replace the old marker text (both runs, including the\n// trailing space run) with the new marker text in one shot, via search.\n// This naturally merges the two runs into one and drops the trailing\n// space, using the first run's formatting for the remaining run.\nconst oldText = \"**ID__AFFARS_pgi_5337_topic_2__ID** \";\nconst newText = \"**ID__AFFARS_AFICC_PGI_5337__ID**\";\nconst matches = body.search(oldText, { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  matches.items[0].insertText(newText, Word.InsertLocation.replace);\n} else {\n  // Fallback: if for some reason the exact text (incl. trailing space)\n  // isn't found (e.g. already edited), just replace the marker text itself.\n  const fallbackMatches = body.search(\"**ID__AFFARS_pgi_5337_topic_2__ID**\", { matchCase: true });\n  fallbackMatches.load(\"items\");\n  await context.sync();\n  if (fallbackMatches.items.length > 0) {\n    fallbackMatches.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// --- 2): update the left indent from 120 twips (6pt) to 225 twips (11.25pt).\nfirstParagraph.leftIndent = 225 / 20;\n\nawait context.sync();\n\n// --- 1) pBdr: Office.js's paragraph.borders only exposes type/color/width\n// (line-style attributes), not the w:space (distance) attribute we need\n// here. Reach the same underlying Word object model the PowerShell/COM\n// side uses (Paragraph.Borders.DistanceFrom{Top,Left,Bottom,Right}) via the\n// native OM bridge that backs this shim, so both automation surfaces\n// produce the identical <w:pBdr><w:top w:space=\"5\"/>...</w:pBdr> markup.\nconst hostHandle = context._root._handle;\nconst anchorJson = JSON.stringify(firstParagraph._anchor);\nconst edges = [\"DistanceFromTop\", \"DistanceFromLeft\", \"DistanceFromBottom\", \"DistanceFromRight\"];\nfor (const edge of edges) {\n  __native.docxOmSet(hostHandle, anchorJson, `Borders.${edge}`, \"5\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (first paragraph of the document, the hidden\n# \"**ID__AFFARS_..._ID**\" marker paragraph):\n#   1. pPr gets a <w:pBdr> with all four edges set to w:space=\"5\"\n#      (no line style / width - just spacing around the paragraph).\n#   2. <w:ind w:left=\"120\"/> -> <w:ind w:left=\"225\"/>  (6pt -> 11.25pt)\n#   3. The two runs \"**ID__AFFARS_pgi_5337_topic_2__ID**\" + \" \" (a\n#      trailing-space-only run) collapse into a single run whose text\n#      becomes \"**ID__AFFARS_AFICC_PGI_5337__ID**\" (no trailing space),\n#      keeping the original run formatting (font/color/size).\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n\n# --- 1) & 3): replace the old marker text (including the trailing space\n# that lives in the second run) with the new marker text (no trailing\n# space) in one Find/Replace over the paragraph's range. This merges the\n# two runs into a single run using the first run's formatting, and drops\n# the separate trailing-space run entirely - matching the target markup.\n$oldText = \"**ID__AFFARS_pgi_5337_topic_2__ID** \"\n$newText = \"**ID__AFFARS_AFICC_PGI_5337__ID**\"\n\n$rng = $p.Range\n$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\nif (-not $found) {\n    # Fallback: exact text (incl. trailing space) not found (e.g. already\n    # edited) - just replace the marker text itself.\n    $rng2 = $p.Range\n    $rng2.Find.Execute(\"**ID__AFFARS_pgi_5337_topic_2__ID**\", $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# Re-fetch paragraph 1 since Find/Replace may shift ranges around.\n$p = $d.Paragraphs(1)\n\n# --- 2): update the left indent from 120 twips (6pt) to 225 twips (11.25pt).\n$p.LeftIndent = 225 / 20\n\n# --- 1) pBdr: set paragraph border spacing on all four edges to 5 twips\n# (<w:top w:space=\"5\"/>, <w:left .../>, <w:bottom .../>, <w:right .../>),\n# with no line style/width set (so only w:space is emitted).\n$p.Borders.DistanceFromTop = 5\n$p.Borders.DistanceFromLeft = 5\n$p.Borders.DistanceFromBottom = 5\n$p.Borders.DistanceFromRight = 5\n"}
